# Refresh the cryptocurrency snapshot (Coin / Link / Price / Volume(1h))
# with the latest scraped values - GitHub Actions scheduled data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new value. "Price" (column D) values are
# text-prefixed with a leading apostrophe so Excel keeps them as literal
# text (e.g. "241.10", "1.000", "6.170", "29.139.97") instead of silently
# re-parsing them as numbers and dropping significant trailing zeros or
# the thousands-grouping dots.
$updates = [ordered]@{
    'D2' = '''29.139.97'
    'E2' = '  +0.70%  '
    'D3' = '''1.831.76'
    'E3' = '  -0.26%  '
    'D4' = '''0.9998'
    'E4' = '  +0.13%  '
    'D5' = '''241.10'
    'E5' = '  -1.71%  '
    'D6' = '''0.6834'
    'E6' = '  -1.65%  '
    'D7' = '''1.001'
    'E7' = '  +0.13%  '
    'D8' = '''0.3016'
    'E8' = '  -0.94%  '
    'D9' = '''0.07476'
    'D10' = '''23.14'
    'E10' = '  -1.07%  '
    'D11' = '''0.07659'
    'E11' = '  -2.12%  '
    'D12' = '''1.835.07'
    'E12' = '  +0.03%  '
    'D13' = '''5.064'
    'E13' = '  -0.85%  '
    'D14' = '''0.6834'
    'E14' = '  +0.23%  '
    'D15' = '''86.94'
    'E15' = '  -6.49%  '
    'D16' = '''6.170'
    'E16' = '  -6.21%  '
    'D17' = '''29.135.18'
    'E17' = '  +0.72%  '
    'D18' = '''0.000008182'
    'E18' = '  -0.97%  '
    'D19' = '''2.079.20'
    'E19' = '  +0.30%  '
    'D20' = '''12.55'
    'E20' = '  -1.21%  '
    'D21' = '''226.74'
    'E21' = '  -6.14%  '
    'E22' = '  +0.11%  '
    'D23' = '''7.439'
    'E23' = '  -0.26%  '
    'D24' = '''1.000'
    'E24' = '  +0.11%  '
    'D25' = '''0.1457'
    'E25' = '  -3.16%  '
    'D26' = '''160.19'
    'E26' = '  +1.43%  '
    'D27' = '''8.759'
    'E27' = '  -0.05%  '
    'D28' = '''18.08'
    'E28' = '  -0.46%  '
    'D29' = '''1.508'
    'E29' = '  -2.12%  '
    'D30' = '''4.268'
    'E30' = '  +1.23%  '
    'D31' = '''4.145'
    'E31' = '  -0.54%  '
    'D32' = '''1.198'
    'E32' = '  +0.44%  '
    'D33' = '''0.05151'
    'E33' = '  +1.03%  '
    'D34' = '''0.7702'
    'E34' = '  -1.48%  '
    'D35' = '''1.842'
    'E35' = '  -0.64%  '
    'D36' = '''1.133'
    'E36' = '  -0.90%  '
    'D37' = '''2.674'
    'E37' = '  -0.82%  '
    'D38' = '''1.311.07'
    'E38' = '  +1.44%  '
    'D39' = '''0.01835'
    'E39' = '  -1.19%  '
    'D40' = '''2.722'
    'E40' = '  +0.76%  '
    'D41' = '''0.9361'
    'E41' = '  -1.54%  '
    'D42' = '''5.788'
    'E42' = '  -5.47%  '
    'D43' = '''104.57'
    'E43' = '  -2.45%  '
    'D44' = '''0.9993'
    'E44' = '  +0.02%  '
    'B45' = 'BabyDogeCoin'
    'C45' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D45' = '''0.00000000123'
    'E45' = '  -0.07%  '
    'B46' = 'Aave'
    'C46' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D46' = '''65.10'
    'E46' = '  +1.87%  '
    'B47' = 'RocketPoolETH'
    'C47' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D47' = '''1.980.69'
    'E47' = '  +0.26%  '
    'B48' = 'Mantle'
    'C48' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D48' = '''0.5198'
    'E48' = '  +0.63%  '
    'D49' = '''9.550'
    'E49' = '  -1.22%  '
    'D50' = '''1.773'
    'E50' = '  +1.14%  '
    'D51' = '''0.05923'
    'E51' = '  +1.04%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
